$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.129.11'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.61%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.281.03'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '154.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +15,378.97%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '304.67'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '94.02'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +0.73%  '
$ws.Range('E8').Value = '  -0.55%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.491'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '34.06'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.22%  '
$ws.Range('E12').Value = '  +0.25%  '
$ws.Range('E13').Value = '  -2.35%  '
$ws.Range('E14').Value = '  -0.34%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.633.17'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.36'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +0.65%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.274.90'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.18%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.791'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +3.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '42.033.42'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.82'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +4.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0918'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +1.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.00'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '68.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.09%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '243.85'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('E25').Value = '  +0.50%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.94'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +1.15%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '24.10'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.76%  '
$ws.Range('E29').Value = '  +5.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.68'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.45%  '
$ws.Range('E31').Value = '  +1.34%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '160.79'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.58%  '
$ws.Range('E33').Value = '  +3.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0752'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('E36').Value = '  +0.77%  '
$ws.Range('E37').Value = '  +3.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.97'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.35%  '
$ws.Range('E39').Value = '  -0.37%  '
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('E41').Value = '  -0.44%  '
$ws.Range('E42').Value = '  +6.78%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '19.86'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +1.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.022.89'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.28'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +11.38%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0284'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.24'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('E48').Value = '  +0.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '53.57'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.28%  '
$ws.Range('E50').Value = '  -1.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.25'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.15%  '
